$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# --- 1) "Pagina login ..." paragraph: drop the proofErr marks / merge
#        its runs into one, then add the new "Gestão de promotores..."
#        list item right after it (same list, numId=3) ---------------
$pLogin = $d.Paragraphs(47)
$xml1 = @"
<w:p $wNs w14:paraId="223345DF" w14:textId="1FFC2B37" w:rsidR="00BC18C9" w:rsidRDefault="00BC18C9" w:rsidP="00BC18C9"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Pagina login - em que só o adm pode adicionar o novo funcionário </w:t></w:r></w:p><w:p $wNs><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Gestão de promotores/funcionários – pôr as opções todas no mesmo form, operações CRUD</w:t></w:r></w:p>
"@
[void]$pLogin.Range.InsertXML($xml1)

# --- 2) "Admin adiciona novo funcionário ..." paragraph: merge runs,
#        drop the proofErr marks (hunk 1 added one extra paragraph, so
#        this paragraph's index shifted from 50 to 51) -----------------
$pAdmin = $d.Paragraphs(51)
$xml2 = @"
<w:p $wNs w14:paraId="6A1C92F3" w14:textId="4FEF4D79" w:rsidR="00721766" w:rsidRDefault="00721766" w:rsidP="00BC18C9"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Admin adiciona novo funcionário (programa devolve palavra passe)</w:t></w:r></w:p>
"@
[void]$pAdmin.Range.InsertXML($xml2)

# --- 3) "Funcionário entra com essa palavra passe ..." paragraph: ----
#        merge runs, drop the proofErr marks (index shifted from 51
#        to 52 for the same reason) ------------------------------------
$pFunc = $d.Paragraphs(52)
$xml3 = @"
<w:p $wNs w14:paraId="683B62B0" w14:textId="2C799F4C" w:rsidR="00721766" w:rsidRDefault="00721766" w:rsidP="00DC0BE3"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Funcionário entra com essa palavra passe, o programa pede para alterar a palavra passe</w:t></w:r></w:p>
"@
[void]$pFunc.Range.InsertXML($xml3)
